$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $s = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $s
}

Set-TextValue $ws "D2" "97.760.40"
Set-TextValue $ws "E2" "  +0.07%  "
Set-TextValue $ws "D3" "3.366.73"
Set-TextValue $ws "E3" "  -1.01%  "
Set-TextValue $ws "E4" "  +0.16%  "
Set-TextValue $ws "D5" "253.31"
Set-TextValue $ws "E5" "  -0.62%  "
Set-TextValue $ws "D6" "660.83"
Set-TextValue $ws "E6" "  +1.21%  "
Set-TextValue $ws "D7" "1.43"
Set-TextValue $ws "E7" "  -2.77%  "
Set-TextValue $ws "D8" "0.429"
Set-TextValue $ws "E8" "  -0.70%  "
Set-TextValue $ws "E9" "  +0.05%  "
Set-TextValue $ws "D10" "1.02"
Set-TextValue $ws "E10" "  -4.49%  "
Set-TextValue $ws "D11" "3.367.19"
Set-TextValue $ws "E11" "  -0.88%  "
Set-TextValue $ws "E12" "  -1.42%  "
Set-TextValue $ws "D13" "42.01"
Set-TextValue $ws "E13" "  +1.00%  "
Set-TextValue $ws "D14" "97.649.89"
Set-TextValue $ws "E14" "  +0.26%  "
Set-TextValue $ws "D15" "6.14"
Set-TextValue $ws "E15" "  -3.20%  "
Set-TextValue $ws "D16" "0.0000257"
Set-TextValue $ws "E16" "  -1.02%  "
Set-TextValue $ws "D17" "3.986.66"
Set-TextValue $ws "E17" "  -1.08%  "
Set-TextValue $ws "D18" "8.81"
Set-TextValue $ws "E18" "  +3.22%  "
Set-TextValue $ws "D19" "3.363.44"
Set-TextValue $ws "E19" "  -1.42%  "
Set-TextValue $ws "D20" "18.01"
Set-TextValue $ws "E20" "  +3.06%  "
Set-TextValue $ws "D21" "0.539"
Set-TextValue $ws "E21" "  +4.66%  "
Set-TextValue $ws "D22" "10.88"
Set-TextValue $ws "E22" "  +1.19%  "
Set-TextValue $ws "D23" "514.14"
Set-TextValue $ws "E23" "  +0.68%  "
Set-TextValue $ws "D24" "3.39"
Set-TextValue $ws "E24" "  -1.79%  "
Set-TextValue $ws "D25" "0.0000202"
Set-TextValue $ws "E25" "  -1.94%  "
Set-TextValue $ws "D26" "6.91"
Set-TextValue $ws "E26" "  +11.39%  "
Set-TextValue $ws "D27" "97.03"
Set-TextValue $ws "E27" "  -2.35%  "
Set-TextValue $ws "D28" "12.44"
Set-TextValue $ws "E28" "  -2.80%  "
Set-TextValue $ws "D29" "3.546.18"
Set-TextValue $ws "E29" "  -0.88%  "
Set-TextValue $ws "D30" "0.147"
Set-TextValue $ws "E30" "  -4.83%  "
Set-TextValue $ws "D31" "11.64"
Set-TextValue $ws "E31" "  +1.67%  "
Set-TextValue $ws "E32" "  +0.70%  "
Set-TextValue $ws "E33" "  -6.67%  "
Set-TextValue $ws "D34" "2.61"
Set-TextValue $ws "E34" "  +14.63%  "
Set-TextValue $ws "B35" "Binance-PegBSC-USD"
Set-TextValue $ws "C35" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws "D35" "1.00"
Set-TextValue $ws "E35" "  +0.29%  "
Set-TextValue $ws "B36" "PolygonEcosystemToken"
Set-TextValue $ws "C36" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws "D36" "0.572"
Set-TextValue $ws "E36" "  +0.71%  "
Set-TextValue $ws "D37" "28.85"
Set-TextValue $ws "E37" "  -2.67%  "
Set-TextValue $ws "D38" "8.01"
Set-TextValue $ws "E38" "  +4.10%  "
Set-TextValue $ws "E39" "  +6.46%  "
Set-TextValue $ws "D40" "523.95"
Set-TextValue $ws "E40" "  -0.19%  "
Set-TextValue $ws "E41" "  -0.38%  "
Set-TextValue $ws "E42" "  +0.08%  "
Set-TextValue $ws "E43" "  +3.98%  "
Set-TextValue $ws "D44" "24.43"
Set-TextValue $ws "E44" "  -1.21%  "
Set-TextValue $ws "D45" "0.862"
Set-TextValue $ws "E45" "  +0.53%  "
Set-TextValue $ws "E46" "  +9.34%  "
Set-TextValue $ws "D47" "5.71"
Set-TextValue $ws "E47" "  +5.31%  "
Set-TextValue $ws "D48" "8.75"
Set-TextValue $ws "E48" "  +6.20%  "
Set-TextValue $ws "D49" "3.64"
Set-TextValue $ws "E49" "  -0.79%  "
Set-TextValue $ws "D50" "53.49"
Set-TextValue $ws "E50" "  +4.51%  "
Set-TextValue $ws "D51" "3.17"
Set-TextValue $ws "E51" "  -3.48%  "
